$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.43 = 30409.67 pesos`n✅ 30409.67 pesos = 7.41 = 967.48 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsTasas.Range("N10").Value = 134.5
$wsTasas.Range("O10").Value = 4090.1
$wsTasas.Range("N12").Value = 4105
$wsTasas.Range("O12").Value = 130.6
